$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "CopperA"

# Small precision corrections to a few existing cells
$ws.Range("J13").Value = 0.9951644108813726
$ws.Range("L13").Value = 0.993158591526912
$ws.Range("O15").Value = 0.9963476568060901

# Append a new data row (row 16) following the same pattern as rows 3-15.
# Copy A15/B15 formatting + shared-string label down into row 16 first,
# then overwrite with the new row's own values.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("B15").Copy($ws.Range("B16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = $ws.Range("B15").Value2

$ws.Range("C16").Value = 1.038744308907424
$ws.Range("D16").Value = 0.9183273690739626
$ws.Range("E16").Value = 1.030356904684638
$ws.Range("F16").Value = 0.9715609114691947
$ws.Range("G16").Value = 1.038744308907424
$ws.Range("H16").Value = 0.9183273690739626
$ws.Range("I16").Value = 1.028838980602037
$ws.Range("J16").Value = 0.981524793711655
$ws.Range("K16").Value = 1.007886970994227
$ws.Range("L16").Value = 0.9385535419865516
$ws.Range("M16").Value = 1.038744308907424
$ws.Range("N16").Value = 0.9743421368793004
$ws.Range("O16").Value = 0.9897473735338049
$ws.Range("P16").Value = 0.9894742226787114
